{"js": "// Change \"VT = 2.24 V\" to \"VT = 2.23 V\" (Plot of HW1Q2 done).\n// The value \"2.24\" only needs a single-digit correction to \"2.23\"; the\n// surrounding run formatting (rFonts/sz/szCs) must stay untouched.\n\nconst body = context.document.body;\n\n// Scope the search to the exact paragraph that holds \"VT = 2.24 V\" so the\n// other \"4\" characters elsewhere in the document (e.g. \"3.45e-4\") are never\n// touched.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"2.24\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the paragraph containing '2.24'.\");\n}\n\nconst results = target.search(\"2.24\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length !== 1) {\n  throw new Error(\"Expected exactly one '2.24' match, found \" + results.items.length);\n}\n\n// Replace just the stale value; formatting is inherited from the run being\n// edited, so rFonts/sz/szCs stay exactly as they were.\nresults.items[0].insertText(\"2.23\", \"Replace\");\nawait context.sync();\n", "ps1": "# Change \"VT = 2.24 V\" to \"VT = 2.23 V\" (Plot of HW1Q2 done).\n# Only the stale value needs to change; run formatting (rFonts/sz/szCs) must\n# stay exactly as-is, so we let Word's normal Find/Replace re-use the run\n# formatting already present at the match.\n\n$d = $word.ActiveDocument\n\n# Find the paragraph that actually contains \"2.24\" so the other \"4\"\n# characters in the document (e.g. \"3.45e-4\") are never touched.\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*2.24*\") {\n        $target = $p.Range\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not find the paragraph containing '2.24'.\"\n}\n\n$find = $target.Find\n$find.ClearFormatting()\n$find.Text = \"2.24\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"2.23\"\n# wdFindContinue = 1, wdReplaceOne = 1\n$find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 1)\n"}
